$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Timestamp column (A2:A97): shift date 2024-06-20 -> 2024-06-25, keep time-of-day ---
$timestamps = @(
    "2024-06-25 00:00:00+00:00",
    "2024-06-25 00:15:00+00:00",
    "2024-06-25 00:30:00+00:00",
    "2024-06-25 00:45:00+00:00",
    "2024-06-25 01:00:00+00:00",
    "2024-06-25 01:15:00+00:00",
    "2024-06-25 01:30:00+00:00",
    "2024-06-25 01:45:00+00:00",
    "2024-06-25 02:00:00+00:00",
    "2024-06-25 02:15:00+00:00",
    "2024-06-25 02:30:00+00:00",
    "2024-06-25 02:45:00+00:00",
    "2024-06-25 03:00:00+00:00",
    "2024-06-25 03:15:00+00:00",
    "2024-06-25 03:30:00+00:00",
    "2024-06-25 03:45:00+00:00",
    "2024-06-25 04:00:00+00:00",
    "2024-06-25 04:15:00+00:00",
    "2024-06-25 04:30:00+00:00",
    "2024-06-25 04:45:00+00:00",
    "2024-06-25 05:00:00+00:00",
    "2024-06-25 05:15:00+00:00",
    "2024-06-25 05:30:00+00:00",
    "2024-06-25 05:45:00+00:00",
    "2024-06-25 06:00:00+00:00",
    "2024-06-25 06:15:00+00:00",
    "2024-06-25 06:30:00+00:00",
    "2024-06-25 06:45:00+00:00",
    "2024-06-25 07:00:00+00:00",
    "2024-06-25 07:15:00+00:00",
    "2024-06-25 07:30:00+00:00",
    "2024-06-25 07:45:00+00:00",
    "2024-06-25 08:00:00+00:00",
    "2024-06-25 08:15:00+00:00",
    "2024-06-25 08:30:00+00:00",
    "2024-06-25 08:45:00+00:00",
    "2024-06-25 09:00:00+00:00",
    "2024-06-25 09:15:00+00:00",
    "2024-06-25 09:30:00+00:00",
    "2024-06-25 09:45:00+00:00",
    "2024-06-25 10:00:00+00:00",
    "2024-06-25 10:15:00+00:00",
    "2024-06-25 10:30:00+00:00",
    "2024-06-25 10:45:00+00:00",
    "2024-06-25 11:00:00+00:00",
    "2024-06-25 11:15:00+00:00",
    "2024-06-25 11:30:00+00:00",
    "2024-06-25 11:45:00+00:00",
    "2024-06-25 12:00:00+00:00",
    "2024-06-25 12:15:00+00:00",
    "2024-06-25 12:30:00+00:00",
    "2024-06-25 12:45:00+00:00",
    "2024-06-25 13:00:00+00:00",
    "2024-06-25 13:15:00+00:00",
    "2024-06-25 13:30:00+00:00",
    "2024-06-25 13:45:00+00:00",
    "2024-06-25 14:00:00+00:00",
    "2024-06-25 14:15:00+00:00",
    "2024-06-25 14:30:00+00:00",
    "2024-06-25 14:45:00+00:00",
    "2024-06-25 15:00:00+00:00",
    "2024-06-25 15:15:00+00:00",
    "2024-06-25 15:30:00+00:00",
    "2024-06-25 15:45:00+00:00",
    "2024-06-25 16:00:00+00:00",
    "2024-06-25 16:15:00+00:00",
    "2024-06-25 16:30:00+00:00",
    "2024-06-25 16:45:00+00:00",
    "2024-06-25 17:00:00+00:00",
    "2024-06-25 17:15:00+00:00",
    "2024-06-25 17:30:00+00:00",
    "2024-06-25 17:45:00+00:00",
    "2024-06-25 18:00:00+00:00",
    "2024-06-25 18:15:00+00:00",
    "2024-06-25 18:30:00+00:00",
    "2024-06-25 18:45:00+00:00",
    "2024-06-25 19:00:00+00:00",
    "2024-06-25 19:15:00+00:00",
    "2024-06-25 19:30:00+00:00",
    "2024-06-25 19:45:00+00:00",
    "2024-06-25 20:00:00+00:00",
    "2024-06-25 20:15:00+00:00",
    "2024-06-25 20:30:00+00:00",
    "2024-06-25 20:45:00+00:00",
    "2024-06-25 21:00:00+00:00",
    "2024-06-25 21:15:00+00:00",
    "2024-06-25 21:30:00+00:00",
    "2024-06-25 21:45:00+00:00",
    "2024-06-25 22:00:00+00:00",
    "2024-06-25 22:15:00+00:00",
    "2024-06-25 22:30:00+00:00",
    "2024-06-25 22:45:00+00:00",
    "2024-06-25 23:00:00+00:00",
    "2024-06-25 23:15:00+00:00",
    "2024-06-25 23:30:00+00:00",
    "2024-06-25 23:45:00+00:00"
)

for ($i = 0; $i -lt $timestamps.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $timestamps[$i]
}


# --- Update numeric data for rows 24-61, columns B:F ---
$rowData = @(
    @(24, 981.592819213867, 0.0009815928192138, 0, 0.0004907964096069, 0.0001226991024017),
    @(25, 8783.786763509113, 0.0087837867635091, 0.0009815928192138, 0.0048826897913614, 0.0012206724478403),
    @(26, 19538.31431070964, 0.0195383143107096, 0.0087837867635091, 0.0141610505371093, 0.0035402626342773),
    @(27, 32863.62298583984, 0.0328636229858398, 0.0195383143107096, 0.0262009686482747, 0.006550242162068599),
    @(28, 91083.99003092447, 0.09108399003092429, 0.0328636229858398, 0.061973806508382, 0.0154934516270954),
    @(29, 148583.052734375, 0.1485830527343749, 0.09108399003092429, 0.1198335213826497, 0.0299583803456624),
    @(30, 211957.6834309896, 0.2119576834309895, 0.1485830527343749, 0.1802703680826823, 0.0450675920206705),
    @(31, 330819.7438151042, 0.3308197438151041, 0.2119576834309895, 0.2713887136230468, 0.0678471784057616),
    @(32, 434977.7942708334, 0.4349777942708333, 0.3308197438151041, 0.3828987690429687, 0.09572469226074209),
    @(33, 704467.9103190105, 0.7044679103190103, 0.4349777942708333, 0.5697228522949218, 0.1424307130737304),
    @(34, 792018.859375, 0.7920188593749999, 0.7044679103190103, 0.7482433848470051, 0.1870608462117512),
    @(35, 971652.407877604, 0.9716524078776041, 0.7920188593749999, 0.8818356336263018, 0.2204589084065754),
    @(36, 1220857.216145833, 1.220857216145833, 0.9716524078776041, 1.096254812011719, 0.2740637030029296),
    @(37, 1435919.516276042, 1.435919516276042, 1.220857216145833, 1.328388366210937, 0.3320970915527343),
    @(38, 1589647.139322917, 1.589647139322917, 1.435919516276042, 1.512783327799479, 0.3781958319498697),
    @(39, 1757917.608723958, 1.757917608723958, 1.589647139322917, 1.673782374023437, 0.4184455935058593),
    @(40, 1924739.938802083, 1.924739938802083, 1.757917608723958, 1.841328773763021, 0.4603321934407552),
    @(41, 2121692.915364583, 2.121692915364583, 1.924739938802083, 2.023216427083333, 0.5058041067708332),
    @(42, 2300498.057291667, 2.300498057291667, 2.121692915364583, 2.211095486328125, 0.5527738715820312),
    @(43, 2222340.139322917, 2.222340139322917, 2.300498057291667, 2.261419098307292, 0.5653547745768228),
    @(44, 2621090.74609375, 2.62109074609375, 2.222340139322917, 2.421715442708333, 0.6054288606770832),
    @(45, 2746057.751302083, 2.746057751302083, 2.62109074609375, 2.683574248697917, 0.670893562174479),
    @(46, 2883473.609375, 2.883473609375, 2.746057751302083, 2.814765680338541, 0.7036914200846354),
    @(47, 2986086.264322917, 2.986086264322917, 2.883473609375, 2.934779936848958, 0.7336949842122394),
    @(48, 3040622.3984375, 3.0406223984375, 2.986086264322917, 3.013354331380208, 0.7533385828450521),
    @(49, 0, 0, 3.0406223984375, 1.52031119921875, 0.3800777998046875),
    @(50, 0, 0, 0, 0, 0),
    @(51, 0, 0, 0, 0, 0),
    @(52, 0, 0, 0, 0, 0),
    @(53, 0, 0, 0, 0, 0),
    @(54, 0, 0, 0, 0, 0),
    @(55, 0, 0, 0, 0, 0),
    @(56, 0, 0, 0, 0, 0),
    @(57, 0, 0, 0, 0, 0),
    @(58, 0, 0, 0, 0, 0),
    @(59, 0, 0, 0, 0, 0),
    @(60, 0, 0, 0, 0, 0),
    @(61, 0, 0, 0, 0, 0)
)

foreach ($rd in $rowData) {
    $r = $rd[0]
    $ws.Cells.Item($r, 2).Value = $rd[1]
    $ws.Cells.Item($r, 3).Value = $rd[2]
    $ws.Cells.Item($r, 4).Value = $rd[3]
    $ws.Cells.Item($r, 5).Value = $rd[4]
    $ws.Cells.Item($r, 6).Value = $rd[5]
}

